$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 183; this shifts the existing rows 183-222
# down to 184-223 (old row 222 becomes new row 223).
$ws.Rows("183").Insert()

# Populate the newly inserted row 183. Its content mirrors the row that
# used to be at 183 (now at 184), except for the date (column D) and the
# volume (column J), which carry the new weekly entry's values.
$ws.Range("A183").Value = 4
$ws.Range("B183").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C183").Value = "Los Lagos"
$ws.Range("D183").Value = 44476
$ws.Range("E183").Value = 10
$ws.Range("F183").Value = 100114013
$ws.Range("G183").Value = "Zanahoria"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 250
$ws.Range("K183").Value = 12000
$ws.Range("L183").Value = 12000
$ws.Range("M183").Value = 12000
$ws.Range("N183").Value = "$/saco 20 kilos"
$ws.Range("O183").Value = "Región de Ñuble"
$ws.Range("P183").Value = 600
$ws.Range("Q183").Value = 20
$ws.Range("R183").Value = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Range("D183").NumberFormat = "YYYY-MM-DD HH:MM:SS"
